# dsa dp and greedy
# Adds two new rows (300. Longest Increasing Subsequence / 763. Partition
# Labels) to the bottom of the LeetCode tracker table on Sheet1, wires up
# their hyperlinks, grows Table2 to cover them, and updates the sheet
# selection the way the authored commit left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 63 - 300. Longest Increasing Subsequence
# ---------------------------------------------------------------------
$ws.Range("A63").Value = "300. Longest Increasing Subsequence"
$ws.Range("B63").Value = "Medium"
$ws.Range("C63").Value = "Dynamic Programming"
$ws.Range("D63").Value = 'Classic DP problem. A subsequence is a sequence that is not necessarily contiguous. Initialize the dp[] with 1 values as default. Maintain a local max var. We perform a nested for loop, outer i from the the right and inner j from the left. If the nums[j] inner value is greater than the outer value nums[i], then consider a new max length for that index dp[i]. Fill dp[i] = Math.max(1+dp[j], dp[i]).'

$link63 = 'https://leetcode.com/problems/longest-increasing-subsequence/solutions/74953/java-solution-dp-simple/ '
$ws.Range("E63").Value = $link63
$ws.Hyperlinks.Add($ws.Range("E63"), $link63, "", "")

# ---------------------------------------------------------------------
# Row 64 - 763. Partition Labels
# ---------------------------------------------------------------------
$ws.Range("A64").Value = "763. Partition Labels"
$ws.Range("B64").Value = "Medium"
$ws.Range("C64").Value = "Greedy"
$ws.Range("D64").Value = 'We care about the last index where each character occurs at, and we can use a HashMap. Do 2 passes: 1. HashMap, 2. Output. We need to update the while loop (it will be a nested while loop) to extend the current partition to cover all the characters that have a last occurrence inside the partition. After i passes j, we find a valid partition and add it to the result list, until i reaches the end of the string.'

$link64 = 'https://leetcode.com/problems/partition-labels/solutions/1868842/java-c-visually-explaineddddd/ '
$ws.Range("E64").Value = $link64
$ws.Hyperlinks.Add($ws.Range("E64"), $link64, "", "")

# ---------------------------------------------------------------------
# Re-apply the row formatting (fills / hyperlink font) from an existing
# "Medium" row onto the two new rows, now that the values + hyperlinks
# are in place (doing this last keeps the hyperlink style index in sync
# with the rest of the table instead of minting a near-duplicate one).
# ---------------------------------------------------------------------
$ws.Range("A2:E2").Copy()
$ws.Range("A63:E64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Grow Table2 ("Question/Difficulty/Pattern/Notes/Link") to include the
# two new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E64"))

# ---------------------------------------------------------------------
# Match the saved selection/viewport from the commit.
# ---------------------------------------------------------------------
$ws.Range("C66").Select()
